$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data: VSCode column header and Python row values ---
$ws.Range("F1").Value = "VSCode"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 4

# --- Styling: header row (rotated, centered, thin border) ---
$rngHeader = $ws.Range("B1:F1")
$rngHeader.Borders.LineStyle = 1
$rngHeader.HorizontalAlignment = -4108
$rngHeader.Orientation = 90

# --- Styling: data grid (centered, thin border) ---
$rngData = $ws.Range("B2:F8")
$rngData.Borders.LineStyle = 1
$rngData.HorizontalAlignment = -4108

# --- Styling: label column (thin border only) ---
$ws.Range("A1:A8").Borders.LineStyle = 1

# --- Column widths for the new grid columns ---
$ws.Range("B1").ColumnWidth = 2.5299479166666665
$ws.Range("C1").ColumnWidth = 3.0729166666666665
$ws.Range("D1:F1").ColumnWidth = 2.5299479166666665

# --- Row height for the (now rotated) header row ---
$ws.Rows.Item(1).RowHeight = 55

# --- Selection as left by the author ---
$ws.Range("E8").Select()
